$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# Update the "Reference" column (A2:A7) with new transaction identifiers
$ws.Range("A2").Value = "Transacción N a915ff59-eb8f-4817-946e-2cb5da370461"
$ws.Range("A3").Value = "Transacción N 1a0ad5f1-0397-4eb0-a8dc-768c9e2ce2a8"
$ws.Range("A4").Value = "Transacción N fbffd5fe-cc28-42ac-ab1b-f794eea3475e"
$ws.Range("A5").Value = "Transacción N fb51aedc-fe4f-4004-aa7b-f2352526267f"
$ws.Range("A6").Value = "Transacción N 33d18ab6-3442-4fe4-ad3a-5f2e129bf3e6"
$ws.Range("A7").Value = "Transacción N 6fdfed77-67a9-4b74-878e-44ed268fc753"

# Update the "UserControlRoom" column (F2:F7) with corrected email address
$ws.Range("F2").Value = "josea.maciast@ecci.edu.co"
$ws.Range("F3").Value = "josea.maciast@ecci.edu.co"
$ws.Range("F4").Value = "josea.maciast@ecci.edu.co"
$ws.Range("F5").Value = "josea.maciast@ecci.edu.co"
$ws.Range("F6").Value = "josea.maciast@ecci.edu.co"
$ws.Range("F7").Value = "josea.maciast@ecci.edu.co"

# Update the selection to reflect the new active cell/selected range
$ws.Activate()
$ws.Range("F2:I10").Select()
